# Fruta / hortaliza, semanal
# Insert a new daily-price record row for "Vega Modelo de Temuco - Uva" at row 464,
# pushing the existing rows 464:556 down to 465:557.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 464 (shifts 464:556 -> 465:557)
$ws.Rows(464).Insert()

# Populate the newly inserted row 464 with the new price record.
$ws.Range("A464").Value = 10
$ws.Range("B464").Value = "Vega Modelo de Temuco"
$ws.Range("C464").Value = "La Araucanía"
$ws.Range("D464").Value = 44522
$ws.Range("E464").Value = 9
$ws.Range("F464").Value = "Fruta"
$ws.Range("G464").Value = 100109
$ws.Range("H464").Value = "Uva"
$ws.Range("I464").Value = 100109001
$ws.Range("J464").Value = "Uva"
$ws.Range("K464").Value = "Superior Seedless"
$ws.Range("L464").Value = "Primera"
$ws.Range("M464").Value = 400
$ws.Range("N464").Value = 27000
$ws.Range("O464").Value = 28000
$ws.Range("P464").Value = 27500
$ws.Range("Q464").Value = "$/bandeja 8 kilos"
$ws.Range("R464").Value = "EE.UU."
$ws.Range("S464").Value = 3438
$ws.Range("T464").Value = 8
